$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value of 45181 for every
# data row (rows 2-215). Update all of them to 45182, matching the diff.
$range = $ws.Range("C2:C215")
$range.Value = 45182
